# Applies "zero_before_threshold" recalculation results to the
# Step3_DataPts_* sheets (0.5 / 0.7 / 0.8 / 0.9 intensity thresholds).
#
# Columns: A=Segment_ID B=Intensity_Threshold C=First_Noticeable_Increase_Index
#          D=Point_Exceeds_Index E=First_Noticeable_Increase_Cumulative_Value
#          F=Point_Exceeds_Cumulative_Value G=Pulse_Width H=Tire_Number ...
#
# Only columns C, E and G change; D (and everything else) is unaffected,
# since Pulse_Width (G) = Point_Exceeds_Index (D) - First_Noticeable_Increase_Index (C).

$wb = $excel.ActiveWorkbook

# Per-row updates that are identical on every Step3_DataPts_* sheet
# (First_Noticeable_Increase_Index / Cumulative_Value don't depend on threshold).
$rowUpdates = @{
    2 = @{ C = 89; E = 0.002956493920201198 }
    3 = @{ C = 87; E = 0.002184565466431852 }
    4 = @{ C = 91; E = 0.01603984421151782 }
    5 = @{ C = 88; E = 0.002205412745566444 }
    6 = @{ C = 89; E = 0.001778565564443984 }
}

# Pulse_Width (G) new values, per sheet name and row.
$pulseWidthUpdates = @{
    "Step3_DataPts_0.5" = @{ 2 = 43; 3 = 39; 4 = 42; 5 = 38; 6 = 42 }
    "Step3_DataPts_0.7" = @{ 2 = 60; 3 = 61; 4 = 59; 5 = 60; 6 = 60 }
    "Step3_DataPts_0.8" = @{ 2 = 69; 3 = 70; 4 = 67; 5 = 66; 6 = 68 }
    "Step3_DataPts_0.9" = @{ 2 = 81; 3 = 82; 4 = 79; 5 = 81; 6 = 80 }
}

foreach ($sheetName in $pulseWidthUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $rowUpdates.Keys) {
        $ws.Range("C$row").Value = $rowUpdates[$row].C
        $ws.Range("E$row").Value = $rowUpdates[$row].E
        $ws.Range("G$row").Value = $pulseWidthUpdates[$sheetName][$row]
    }
}
